# Update odds/statistics values in the weekly FlashScore games export.
# Applies the numeric corrections from rows 2-4 as described by the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    # Row 2
    "G2"  = 3.1
    "H2"  = 3.25
    "I2"  = 2.3
    "J2"  = 3.6
    "L2"  = 3
    "M2"  = 1.07
    "N2"  = 7.5
    "AB2" = 29
    "AG2" = 15
    "AJ2" = 7.5
    "AK2" = 11
    "AL2" = 10
    "AN2" = 21

    # Row 3
    "Q3"  = 2.38
    "R3"  = 1.53
    "S3"  = 4.5
    "T3"  = 1.17
    "AP3" = 1.85
    "AQ3" = 2

    # Row 4
    "H4"  = 3.1
    "I4"  = 2.8
    "K4"  = 2.2
    "L4"  = 3.4
    "M4"  = 1.05
    "N4"  = 11
    "O4"  = 1.25
    "P4"  = 3.75
    "Q4"  = 1.88
    "R4"  = 1.98
    "S4"  = 3
    "T4"  = 1.36
    "U4"  = 1.36
    "V4"  = 3
    "W4"  = 1.62
    "X4"  = 2.2
    "Y4"  = 10
    "AD4" = 26
    "AE4" = 11
    "AF4" = 6
    "AI4" = 151
    "AM4" = 29
    "AO4" = 29
}

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}
